# 20241121_Sink_Isolates_AR_22-6.xlsx — antibiotic-resistance re-read update.
#
# The lab re-measured three isolate/antibiotic combos ("Redo (Nov22)": isolate
# 22-6-a vs Carb, 22-6-e vs Carb, 22-6-f vs Kan) on the "Sheet1" tab, logged
# the new triplicate readings in a small side-table (columns Q:U, rows 29-35),
# and then propagated the corrected averages back into the main data block
# (columns B:D) on both "Sheet1" and "Sheet2". "Sheet3" only has its view
# state touched.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# Sheet1: corrected triplicate values for isolate 22-6-b (row 31, Carb) and
# for row 33 (Kan) pulled in from the redo.
# ---------------------------------------------------------------------------
$ws1.Range("B31").Value = 1.36650002002716
$ws1.Range("C31").Value = 1.31319999694824
$ws1.Range("D31").Value = 1.2639000415802

$ws1.Range("E33").Value = 0.0507999993860722
$ws1.Range("F33").Value = 0.050500001758337
$ws1.Range("G33").Value = 0.0516000017523766

# Side table recording the Nov-22 redo readings ------------------------------

# Header label for the redo block
$ws1.Range("S29").Value = "Redo (Nov22)"

# Replicate numbers header (styled like the other header rows, e.g. A29:M29)
$ws1.Range("S30").Value = 10
$ws1.Range("T30").Value = 11
$ws1.Range("U30").Value = 12
$ws1.Range("S30:U30").Interior.Color = 8421504
$ws1.Range("S30:U30").Font.Color = 16777215

# LB control row
$ws1.Range("R31").Value = "LB"
$ws1.Range("S31").Value = 0.045099999755620956
$ws1.Range("T31").Value = 0.04399999976158142
$ws1.Range("U31").Value = 0.045099999755620956

# 22-6-a vs Carb redo
$ws1.Range("Q32").Value = "22-6-a"
$ws1.Range("R32").Value = "Carb"
$ws1.Range("S32").Value = 1.2532999515533401
$ws1.Range("T32").Value = 1.33109998703003
$ws1.Range("U32").Value = 1.24049997329712

# 22-6-e vs Carb redo
$ws1.Range("Q33").Value = "22-6-e"
$ws1.Range("R33").Value = "Carb"
$ws1.Range("S33").Value = 1.36650002002716
$ws1.Range("T33").Value = 1.31319999694824
$ws1.Range("U33").Value = 1.2639000415802

# 22-6-f vs Kan redo
$ws1.Range("Q34").Value = "22-6-f"
$ws1.Range("R34").Value = "Kan"
$ws1.Range("S34").Value = 0.0507999993860722
$ws1.Range("T34").Value = 0.050500001758337
$ws1.Range("U34").Value = 0.0516000017523766

# LB control row
$ws1.Range("R35").Value = "LB"
$ws1.Range("S35").Value = 0.04439999908208847
$ws1.Range("T35").Value = 0.043800000101327896
$ws1.Range("U35").Value = 0.04439999908208847

# ---------------------------------------------------------------------------
# Sheet2: same corrected triplicate (isolate 22-6-a vs Carb) propagated into
# the main block, and the existing redo columns re-saved at full precision.
# ---------------------------------------------------------------------------
$ws2.Range("B31").Value = 1.2532999515533401
$ws2.Range("C31").Value = 1.33109998703003
$ws2.Range("D31").Value = 1.24049997329712

$ws2.Range("R31").Value = 1.2532999515533401
$ws2.Range("S31").Value = 1.33109998703003
$ws2.Range("T31").Value = 1.24049997329712

# ---------------------------------------------------------------------------
# View state: the redo work was done with "Sheet1" as the active/visible tab.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("P28:T35").Select()

$ws3.Activate()
$ws3.Range("K49").Select()

$ws1.Activate()
$ws1.Range("I47").Select()
